# Adding the changes we made on may 9th
# Replace the sensor data block with the May-9th dataset: 7 new rows
# inserted at the top (pushing the original 20 rows down) plus 3 new
# rows appended at the end, for 30 data rows total (was 20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 30,8
$arr[0,0] = 0
$arr[0,1] = "falling"
$arr[0,2] = -0.7812347412109375
$arr[0,3] = 2.500675392150879
$arr[0,4] = 2.085890746116638
$arr[0,5] = -0.0421497002243995
$arr[0,6] = 0.0813977941870689
$arr[0,7] = -0.0167987942695617
$arr[1,0] = 100
$arr[1,1] = "falling"
$arr[1,2] = -0.6001354455947877
$arr[1,3] = 2.425701588392258
$arr[1,4] = 1.942269176244736
$arr[1,5] = -0.0009162978967651
$arr[1,6] = -0.0372627787292003
$arr[1,7] = -0.051312681287527
$arr[2,0] = 200
$arr[2,1] = "falling"
$arr[2,2] = -0.4221334457397461
$arr[2,3] = 2.401677787303925
$arr[2,4] = 1.923809313774109
$arr[2,5] = -0.0261144898831844
$arr[2,6] = 0.0638354197144508
$arr[2,7] = -0.0560468845069408
$arr[3,0] = 300
$arr[3,1] = "falling"
$arr[3,2] = -0.4831114292144776
$arr[3,3] = 2.420794081687927
$arr[3,4] = 2.053104478120804
$arr[3,5] = 0.0224492978304624
$arr[3,6] = 0.0224492978304624
$arr[3,7] = 0.0267253536731004
$arr[4,0] = 400
$arr[4,1] = "falling"
$arr[4,2] = -0.4990121841430664
$arr[4,3] = 2.45754919052124
$arr[4,4] = 2.050947427749634
$arr[4,5] = 0.0284052342176437
$arr[4,6] = -0.0120645882561802
$arr[4,7] = 0.0143553335219621
$arr[5,0] = 500
$arr[5,1] = "falling"
$arr[5,2] = -0.4582573175430297
$arr[5,3] = 2.450900214910507
$arr[5,4] = 1.956571793556213
$arr[5,5] = 0.0213802829384803
$arr[5,6] = 0.0073303831741213
$arr[5,7] = -0.0117591563612222
$arr[6,0] = 600
$arr[6,1] = "falling"
$arr[6,2] = -0.5409791469573973
$arr[6,3] = 2.33053719997406
$arr[6,4] = 1.996131032705307
$arr[6,5] = -0.0149661982432007
$arr[6,6] = 0.042302418500185
$arr[6,7] = 0.0258090570569038
$arr[7,0] = 700
$arr[7,1] = "falling"
$arr[7,2] = -0.6303038597106934
$arr[7,3] = 2.39460033774376
$arr[7,4] = 2.063519307971001
$arr[7,5] = 0.0442877300083637
$arr[7,6] = 0.0798706337809562
$arr[7,7] = 0.0123700210824608
$arr[8,0] = 800
$arr[8,1] = "falling"
$arr[8,2] = -0.5232744216918948
$arr[8,3] = 2.33260555267334
$arr[8,4] = 1.824701523780822
$arr[8,5] = 0.0595593601465225
$arr[8,6] = 0.0303905457258224
$arr[8,7] = 0.0181732401251792
$arr[9,0] = 900
$arr[9,1] = "falling"
$arr[9,2] = -0.7410809755325318
$arr[9,3] = 2.134066888689995
$arr[9,4] = 1.285716485977173
$arr[9,5] = -0.00534507073462
$arr[9,6] = -0.0339030213654041
$arr[9,7] = 0.0303905457258224
$arr[10,0] = 1000
$arr[10,1] = "falling"
$arr[10,2] = -0.2835536479949948
$arr[10,3] = 1.980954867601395
$arr[10,4] = 0.2343389749526972
$arr[10,5] = 0.0459676086902618
$arr[10,6] = -0.0532979927957057
$arr[10,7] = 0.1757764816284179
$arr[11,0] = 1100
$arr[11,1] = "falling"
$arr[11,2] = -1.529338419437408
$arr[11,3] = 2.245944246649742
$arr[11,4] = 0.1200629770755768
$arr[11,5] = -0.1426370292901992
$arr[11,6] = -0.4100432991981506
$arr[11,7] = -0.1085812970995903
$arr[12,0] = 1200
$arr[12,1] = "falling"
$arr[12,2] = -1.762292957305908
$arr[12,3] = 1.919589138031006
$arr[12,4] = -0.2954926252365113
$arr[12,5] = -0.1237002089619636
$arr[12,6] = -0.3216205537319183
$arr[12,7] = -0.0311541277915239
$arr[13,0] = 1300
$arr[13,1] = "falling"
$arr[13,2] = -3.197947156429292
$arr[13,3] = 1.681933805346489
$arr[13,4] = -0.3819067515432834
$arr[13,5] = -0.2638937830924988
$arr[13,6] = -0.4506658315658569
$arr[13,7] = -0.1346957832574844
$arr[14,0] = 1400
$arr[14,1] = "falling"
$arr[14,2] = -4.11650104522705
$arr[14,3] = 1.521171301603317
$arr[14,4] = -0.9332275912165637
$arr[14,5] = -0.3332269787788391
$arr[14,6] = -0.4558582007884979
$arr[14,7] = -0.0723875313997268
$arr[15,0] = 1500
$arr[15,1] = "falling"
$arr[15,2] = -4.768530690670014
$arr[15,3] = 1.473171654343605
$arr[15,4] = -0.2953257039189335
$arr[15,5] = 0.0320704244077205
$arr[15,6] = -0.7147123217582703
$arr[15,7] = 0.0733038261532783
$arr[16,0] = 1600
$arr[16,1] = "falling"
$arr[16,2] = -2.971992969512948
$arr[16,3] = 1.761710226535796
$arr[16,4] = -0.3112654983997316
$arr[16,5] = 0.2872593700885772
$arr[16,6] = 0.7376197576522827
$arr[16,7] = -0.3023782968521118
$arr[17,0] = 1700
$arr[17,1] = "falling"
$arr[17,2] = -2.837122094631188
$arr[17,3] = 1.392113929986956
$arr[17,4] = 2.611395421624172
$arr[17,5] = -1.120784997940064
$arr[17,6] = 0.9436340928077698
$arr[17,7] = -2.777909755706787
$arr[18,0] = 1800
$arr[18,1] = "falling"
$arr[18,2] = 0.09898402690889618
$arr[18,3] = 2.09258412122727
$arr[18,4] = 5.970039045810694
$arr[18,5] = -1.083369493484497
$arr[18,6] = 2.2501220703125
$arr[18,7] = -0.494189977645874
$arr[19,0] = 1900
$arr[19,1] = "falling"
$arr[19,2] = 9.021431350707996
$arr[19,3] = 3.897843426465977
$arr[19,4] = 2.894297271966933
$arr[19,5] = -0.2246456891298294
$arr[19,6] = 0.6049093008041382
$arr[19,7] = -0.3110831379890442
$arr[20,0] = 2000
$arr[20,1] = "falling"
$arr[20,2] = 3.752513694763174
$arr[20,3] = 1.26214294433594
$arr[20,4] = 2.215038943290709
$arr[20,5] = 0.0445931628346443
$arr[20,6] = -0.2460259795188903
$arr[20,7] = 0.2205223590135574
$arr[21,0] = 2100
$arr[21,1] = "falling"
$arr[21,2] = -0.8335286378860537
$arr[21,3] = 2.680000334978103
$arr[21,4] = 1.48247443139553
$arr[21,5] = 0.0360410511493682
$arr[21,6] = 1.270752429962158
$arr[21,7] = 0.1033889427781105
$arr[22,0] = 2200
$arr[22,1] = "falling"
$arr[22,2] = 1.149217176437378
$arr[22,3] = 2.79909211397171
$arr[22,4] = 1.387107414007188
$arr[22,5] = 0.1259909570217132
$arr[22,6] = 0.7915286421775818
$arr[22,7] = 0.2243402600288391
$arr[23,0] = 2300
$arr[23,1] = "falling"
$arr[23,2] = 2.04998896121979
$arr[23,3] = 3.530032843351366
$arr[23,4] = -0.8531217783689553
$arr[23,5] = 0.180816113948822
$arr[23,6] = -0.5984952449798584
$arr[23,7] = -0.3129157125949859
$arr[24,0] = 2400
$arr[24,1] = "falling"
$arr[24,2] = -0.119726562500003
$arr[24,3] = 1.84428286552429
$arr[24,4] = 0.8842907547950768
$arr[24,5] = 0.1476766765117645
$arr[24,6] = 0.3284927904605865
$arr[24,7] = 0.2247984111309051
$arr[25,0] = 2500
$arr[25,1] = "falling"
$arr[25,2] = 0.4374212741851807
$arr[25,3] = 2.572291845083237
$arr[25,4] = 0.7996354326605797
$arr[25,5] = 0.0265726372599601
$arr[25,6] = 0.0389426611363887
$arr[25,7] = 0.1386664062738418
$arr[26,0] = 2600
$arr[26,1] = "falling"
$arr[26,2] = 0.3087260723114021
$arr[26,3] = 2.726021051406861
$arr[26,4] = 0.5522446408867842
$arr[26,5] = -0.0685696229338646
$arr[26,6] = -0.0618501044809818
$arr[26,7] = -0.1456913650035858
$arr[27,0] = 2700
$arr[27,1] = "falling"
$arr[27,2] = 0.1567803621292113
$arr[27,3] = 2.781254351139069
$arr[27,4] = 0.9996474064886576
$arr[27,5] = -0.050854530185461
$arr[27,6] = -0.0387899428606033
$arr[27,7] = -0.0474947728216648
$arr[28,0] = 2800
$arr[28,1] = "falling"
$arr[28,2] = 0.246018409729004
$arr[28,3] = 2.646198272705077
$arr[28,4] = 1.203094172477722
$arr[28,5] = -0.0007635815418325
$arr[28,6] = 0.0126754539087414
$arr[28,7] = 0.0435241498053073
$arr[29,0] = 2900
$arr[29,1] = "falling"
$arr[29,2] = 0.4407022714614867
$arr[29,3] = 2.732884711027145
$arr[29,4] = 1.124532252550125
$arr[29,5] = 0.027030786499381
$arr[29,6] = 0.0493273697793483
$arr[29,7] = -0.0366519130766391
$ws.Range("A2:H31").Value = $arr
